$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $style = $range.Style
    $range.Value = "'" + $value
    $range.Style = $style
}

# Row 2
Set-TextValue $ws.Range("D2") "310.29"
Set-TextValue $ws.Range("E2") "0.49%"
# Row 3
Set-TextValue $ws.Range("D3") "41.15"
Set-TextValue $ws.Range("E3") "-0.18%"
# Row 4
Set-TextValue $ws.Range("D4") "5.215"
Set-TextValue $ws.Range("E4") "1.80%"
# Row 5
Set-TextValue $ws.Range("D5") "0.07684"
Set-TextValue $ws.Range("E5") "0.60%"
# Row 6
Set-TextValue $ws.Range("B6") "GateToken"
Set-TextValue $ws.Range("C6") "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
Set-TextValue $ws.Range("D6") "4.296"
Set-TextValue $ws.Range("E6") "0.57%"
# Row 7
Set-TextValue $ws.Range("B7") "FTXToken"
Set-TextValue $ws.Range("C7") "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
Set-TextValue $ws.Range("D7") "1.745"
Set-TextValue $ws.Range("E7") "7.75%"
# Row 8
Set-TextValue $ws.Range("B8") "MXToken"
Set-TextValue $ws.Range("C8") "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
Set-TextValue $ws.Range("D8") "0.9274"
Set-TextValue $ws.Range("E8") "2.32%"
# Row 9
Set-TextValue $ws.Range("B9") "BTSEToken"
Set-TextValue $ws.Range("C9") "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
Set-TextValue $ws.Range("D9") "2.425"
Set-TextValue $ws.Range("E9") "-1.22%"
# Row 10
Set-TextValue $ws.Range("B10") "LiechtensteinCryptoassetsExchange"
Set-TextValue $ws.Range("C10") "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
Set-TextValue $ws.Range("D10") "0.1276"
Set-TextValue $ws.Range("E10") "13.06%"
# Row 11
Set-TextValue $ws.Range("B11") "WazirX"
Set-TextValue $ws.Range("C11") "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
Set-TextValue $ws.Range("D11") "0.1824"
Set-TextValue $ws.Range("E11") "1.30%"
# Row 12
Set-TextValue $ws.Range("B12") "MandalaExchangeToken"
Set-TextValue $ws.Range("C12") "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
Set-TextValue $ws.Range("D12") "0.09111"
Set-TextValue $ws.Range("E12") "-0.16%"
# Row 13
Set-TextValue $ws.Range("B13") "BitrueCoin"
Set-TextValue $ws.Range("C13") "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
Set-TextValue $ws.Range("D13") "0.04167"
Set-TextValue $ws.Range("E13") "-1.92%"
# Row 14
Set-TextValue $ws.Range("B14") "BitMartToken"
Set-TextValue $ws.Range("C14") "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
Set-TextValue $ws.Range("D14") "0.1051"
Set-TextValue $ws.Range("E14") "0.18%"
# Row 15
Set-TextValue $ws.Range("B15") "BitForexToken"
Set-TextValue $ws.Range("C15") "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
Set-TextValue $ws.Range("D15") "0.001291"
Set-TextValue $ws.Range("E15") "3.07%"
# Row 16
Set-TextValue $ws.Range("B16") "TigerCash"
Set-TextValue $ws.Range("C16") "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
Set-TextValue $ws.Range("D16") "0.005888"
Set-TextValue $ws.Range("E16") "2.01%"
# Row 17
Set-TextValue $ws.Range("B17") "LEO"
Set-TextValue $ws.Range("C17") "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
Set-TextValue $ws.Range("D17") "3.353"
Set-TextValue $ws.Range("E17") "0.20%"
# Row 18
Set-TextValue $ws.Range("E18") "-0.38%"
# Row 19
Set-TextValue $ws.Range("D19") "7.392"
Set-TextValue $ws.Range("E19") "10.06%"
# Row 20
Set-TextValue $ws.Range("E20") "-0.86%"
# Row 21
Set-TextValue $ws.Range("E21") "-0.52%"
# Row 22
Set-TextValue $ws.Range("D22") "0.04012"
Set-TextValue $ws.Range("E22") "-1.15%"
# Row 23
Set-TextValue $ws.Range("D23") "0.001266"
Set-TextValue $ws.Range("E23") "-0.01%"
# Row 24
Set-TextValue $ws.Range("D24") "0.004096"
Set-TextValue $ws.Range("E24") "1.35%"
# Row 25
Set-TextValue $ws.Range("D25") "0.0001270"
Set-TextValue $ws.Range("E25") "-0.11%"
# Row 38
Set-TextValue $ws.Range("D38") "0.02527"
Set-TextValue $ws.Range("E38") "4.20%"
# Row 39
Set-TextValue $ws.Range("D39") "0.05316"
Set-TextValue $ws.Range("E39") "1.48%"
# Row 40
Set-TextValue $ws.Range("D40") "0.007844"
Set-TextValue $ws.Range("E40") "0.57%"
# Row 41
Set-TextValue $ws.Range("D41") "0.1311"
Set-TextValue $ws.Range("E41") "0.79%"
# Row 42
Set-TextValue $ws.Range("D42") "0.006637"
Set-TextValue $ws.Range("E42") "1.55%"
# Row 43
Set-TextValue $ws.Range("D43") "0.002051"
Set-TextValue $ws.Range("E43") "5.10%"
# Row 44
Set-TextValue $ws.Range("D44") "0.008096"
Set-TextValue $ws.Range("E44") "6.71%"
# Row 45
Set-TextValue $ws.Range("D45") "0.3085"
Set-TextValue $ws.Range("E45") "-0.06%"
# Row 46
Set-TextValue $ws.Range("D46") "0.00006779"
Set-TextValue $ws.Range("E46") "0.06%"
# Row 47
Set-TextValue $ws.Range("D47") "0.00000000750"
Set-TextValue $ws.Range("E47") "-0.11%"
# Row 48
Set-TextValue $ws.Range("D48") "0.2278"
Set-TextValue $ws.Range("E48") "240.49%"
# Row 50
Set-TextValue $ws.Range("D50") "0.00002099"
Set-TextValue $ws.Range("E50") "-0.11%"
# Row 51
Set-TextValue $ws.Range("D51") "0.0001999"
Set-TextValue $ws.Range("E51") "-0.11%"
